$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 187; this shifts the existing rows 187-237
# down to 188-238 and extends the used range to R238, matching the
# weekly data refresh described in the commit message.
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with this week's record.
$ws.Range("A187").Value = 9
$ws.Range("B187").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C187").Value = "Metropolitana"
$ws.Range("D187").Value = 44508
$ws.Range("E187").Value = 13
$ws.Range("F187").Value = 100112044
$ws.Range("G187").Value = "Perejil"
$ws.Range("H187").Value = "Sin especificar"
$ws.Range("I187").Value = "Primera"
$ws.Range("J187").Value = 79
$ws.Range("K187").Value = 9000
$ws.Range("L187").Value = 10000
$ws.Range("M187").Value = 9430
$ws.Range("N187").Value = "$/docena de atados"
$ws.Range("O187").Value = "Región Metropolitana"
$ws.Range("P187").Value = 3143
$ws.Range("Q187").Value = 3
$ws.Range("R187").Value = "Hortaliza"

# Preserve the date-formatted style used throughout column D.
$ws.Range("D187").NumberFormat = "YYYY-MM-DD HH:MM:SS"
